$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("R4").Value = 0.057927248158369672
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Font.Name = "Times New Roman"
$ws.Range("R4").Font.Size = 9
$ws.Range("R4").Font.Bold = $true
$ws.Range("R4").Font.Bold = $false
